# =====================================================================
# Edit script: updates "杭州-漫展信息" workbook per commit diff.
#   Sheet 1 (展览)     -> update "想去人数" (F) counters
#   Sheet 2 (演出)     -> remove a duplicated cancelled listing, add a
#                         newly scraped listing (2024CJMF·不止音乐节),
#                         and refresh counters/links for the rest
#   Sheet 3 (本地生活) -> update "想去人数" (F) counters
#   Sheet 4 (全部类型) -> update "想去人数" (F) counters (mirrors 1-3)
# =====================================================================

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) - numeric counter refresh only
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$sheet1Updates = @{
    "F2"  = 12885
    "F3"  = 7212
    "F10" = 1014
    "F11" = 151
    "F12" = 360
    "F13" = 1026
    "F14" = 8
    "F16" = 1023
    "F17" = 511
    "F18" = 256
    "F19" = 374
    "F22" = 315
    "F24" = 197
    "F25" = 380
    "F26" = 5261
    "F27" = 73
    "F28" = 1445
    "F29" = 316
    "F30" = 1714
    "F31" = 76
    "F32" = 70
    "F33" = 1374
    "F36" = 601
    "F38" = 3743
}
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

# -----------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# The duplicated "Marcin Patrzalek" row (rows 2 and 3 are identical
# duplicates of a listing that no longer appears) is removed; the
# table shifts up by two rows.
$ws2.Range("A2:I2").EntireRow.Delete()
$ws2.Range("A2:I2").EntireRow.Delete()

# A newly scraped listing ("2024CJMF·不止音乐节") is inserted between
# the "红楼梦" listing and the "吉卜力" listing (new row 7).
$ws2.Range("A7:I7").EntireRow.Insert()

# Fill the brand-new row with its data. B7 holds a plain "yyyy-mm-dd"
# looking string, so force text formatting first or Excel will coerce
# it into a date serial number.
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "2024-09-15"
$ws2.Range("C7").Value = "杭州·2024CJMF·不止音乐节"
$ws2.Range("D7").Value = "塘子堰路177号 华数产业园隔壁大草坪"
$ws2.Range("E7").Value = "2024.09.15 13:00-09.16 21:40"
$ws2.Range("F7").Value = 61
$ws2.Range("G7").Value = "不可售"
$ws2.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=90522"
$ws2.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202408/3PmG2Bq51723192884141.jpeg"

# Fix up the formatting of the new index cell (A7) to match the rest
# of the index column (bold + bordered), which the row-insert only
# partially carried over.
$ws2.Range("A6").Copy()
$ws2.Range("A7").PasteSpecial(-4122)

# Renumber the index column (A) sequentially for all 19 data rows,
# since deleting/inserting rows does not renumber the literal values
# that were already stored in them.
for ($i = 2; $i -le 20; $i++) {
    $ws2.Range("A$i").Value = $i - 1
}

# Refresh the "想去人数" (F) counter for the "世界经典原版音乐剧《猫》
# CATS" listing, now sitting at row 18 after the shift above.
$ws2.Range("F18").Value = 22

# -----------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) - numeric counter refresh only
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 9292
$ws3.Range("F3").Value = 562
$ws3.Range("F4").Value = 2040

# -----------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - numeric counter refresh only
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$sheet4Updates = @{
    "F2"  = 9292
    "F3"  = 562
    "F4"  = 2040
    "F5"  = 12885
    "F6"  = 7212
    "F10" = 1014
    "F11" = 151
    "F12" = 360
    "F13" = 1026
    "F14" = 8
    "F16" = 1023
    "F17" = 256
    "F18" = 374
    "F21" = 315
    "F26" = 197
    "F27" = 380
    "F28" = 5261
    "F29" = 73
    "F30" = 1445
    "F33" = 316
    "F35" = 1714
    "F36" = 76
    "F37" = 70
    "F38" = 1374
    "F40" = 601
    "F47" = 3743
    "F48" = 22
}
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
